# The underlying source data got re-synced: within each of the two blocks of
# rows (27-30 and 45-47) the records rotate by one position - the first
# row's record moves to the end of its block, and the others shift up.
#
#   new(27) = old(28)   new(45) = old(46)
#   new(28) = old(29)   new(46) = old(47)
#   new(29) = old(30)   new(47) = old(45)
#   new(30) = old(27)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: rows 27-30 ---------------------------------------------------

# Row 27 <= old row 28
$ws.Range("A27").Value = 130757159
$ws.Range("B27").Value = 79243
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = "Garnlav"
$ws.Range("G27").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Ach.) Ach."
$ws.Range("M27").Value = ""
$ws.Range("P27").Value = "Kråkbackarna, Dlr"
$ws.Range("Q27").Value = 490482
$ws.Range("R27").Value = 6763574
$ws.Range("Z27").Value = "15:01"
$ws.Range("AB27").Value = "15:01"
$ws.Range("AC27").Value = ""

# Row 28 <= old row 29
$ws.Range("A28").Value = 130789472
$ws.Range("Q28").Value = 490450
$ws.Range("R28").Value = 6763926
$ws.Range("Z28").Value = "14:07"
$ws.Range("AB28").Value = "14:07"
$ws.Range("AF28").Value = ""
$ws.Range("AW28").Value = "Bo karlstens"
$ws.Range("AX28").Value = "Bo karlstens, Håkan Thenander"

# Row 29 <= old row 30
$ws.Range("A29").Value = 130754796
$ws.Range("Q29").Value = 490463
$ws.Range("R29").Value = 6763939
$ws.Range("Z29").Value = "11:43"
$ws.Range("AB29").Value = "11:43"
$ws.Range("AF29").ClearContents()
$ws.Range("AW29").Value = "Håkan Thenander"
$ws.Range("AX29").Value = "Håkan Thenander, Bo karlstens"

# Row 30 <= old row 27
$ws.Range("A30").Value = 130755667
$ws.Range("B30").Value = 57884
$ws.Range("E30").Value = 100109
$ws.Range("F30").Value = "Tretåig hackspett"
$ws.Range("G30").Value = "Picoides tridactylus"
$ws.Range("H30").Value = "(Linnaeus, 1758)"
$ws.Range("M30").Value = "färska spår"
$ws.Range("P30").Value = "Prikattmyren, Dlr"
$ws.Range("Q30").Value = 490444
$ws.Range("R30").Value = 6763770
$ws.Range("AC30").Value = "2 bilder"

# --- Block 2: rows 45-47 ----------------------------------------------------

# Row 45 <= old row 46
$ws.Range("A45").Value = 130754014
$ws.Range("B45").Value = 79243
$ws.Range("D45").Value = "NT"
$ws.Range("E45").Value = 6425
$ws.Range("F45").Value = "Garnlav"
$ws.Range("G45").Value = "Alectoria sarmentosa"
$ws.Range("H45").Value = "(Ach.) Ach."
$ws.Range("M45").Value = ""
$ws.Range("P45").Value = "Kråkbackarna, Dlr"
$ws.Range("Q45").Value = 490548
$ws.Range("R45").Value = 6763654
$ws.Range("Z45").Value = "11:43"
$ws.Range("AB45").Value = "11:43"

# Row 46 <= old row 47
$ws.Range("A46").Value = 130754953
$ws.Range("P46").Value = "Prikattmyren, Dlr"
$ws.Range("Q46").Value = 490440
$ws.Range("R46").Value = 6764028
$ws.Range("AC46").Value = "1 bild. Rikligt på gran"

# Row 47 <= old row 45
$ws.Range("A47").Value = 130758028
$ws.Range("B47").Value = 8451
$ws.Range("D47").Value = "LC"
$ws.Range("E47").Value = 106545
$ws.Range("F47").Value = "Mindre märgborre"
$ws.Range("G47").Value = "Tomicus minor"
$ws.Range("H47").Value = "(Hartig, 1834)"
$ws.Range("M47").Value = "äldre gnagspår"
$ws.Range("P47").Value = "Brunnvasselänget, Dlr"
$ws.Range("Q47").Value = 490175
$ws.Range("R47").Value = 6763613
$ws.Range("Z47").Value = "15:01"
$ws.Range("AB47").Value = "15:01"
$ws.Range("AC47").ClearContents()
